$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (shifts existing rows 5..73 down to 6..74)
$ws.Rows.Item(5).Insert()

# New journal name mapping: raw PubMed-style title in column A,
# clean formatted title in column B.
$ws.Cells.Item(5, 2).Value = "Archives of Otolaryngology Head and Neck Surgery"
$ws.Cells.Item(5, 1).Value = "Archives of otolaryngology--head & neck surgery"

# Re-apply the existing sort (keyed on column A) so the sheet's persisted
# sort state keeps tracking the full A2:B74 data range after the insert.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2"))
$ws.Sort.SetRange($ws.Range("A2:B74"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

$ws.Range("A6").Select()
